# Rename the "Specialist" ecology category to "Strict alpine".
# The workbook's "species24" sheet has column G = "ecology" with values
# "Specialist" or "Generalist"; every "Specialist" becomes "Strict alpine".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("species24")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Specialist") {
        $cell.Value2 = "Strict alpine"
    }
}

# Reflect the selection change recorded for the sheet: the author ended up
# with the whole ecology column (G) selected.
$ws.Range("G:G").Select()
